$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B/C (Coin name / Link) updates ---
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "68.771.30"
$ws.Range("D3").Value = "3.817.57"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.61"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.83"
$ws.Range("D7").Value = "3.814.79"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.43"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000262"
$ws.Range("D15").Value = "4.457.48"
$ws.Range("D16").Value = "3.798.39"
$ws.Range("D17").Value = "68.718.74"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.96"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.79"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "466.75"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.703"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000156"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.37"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.18"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.97"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.15"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.34"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "30.15"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.19"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.19"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D36").Value = "3.767.36"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.53"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.138"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.82"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.67"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.93"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.92"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.45"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "147.11"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "392.31"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000267"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("E7").Value = "  +1.12%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("E13").Value = "  -2.36%  "
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("E18").Value = "  -1.80%  "
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("E24").Value = "  +8.45%  "
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("E26").Value = "  -3.36%  "
$ws.Range("E27").Value = "  -2.06%  "
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("E31").Value = "  -1.78%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  -3.20%  "
$ws.Range("E34").Value = "  -1.63%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("E36").Value = "  +1.01%  "
$ws.Range("E37").Value = "  -1.49%  "
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("E45").Value = "  +15.25%  "
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("E47").Value = "  +2.99%  "
$ws.Range("E48").Value = "  -2.10%  "
$ws.Range("E49").Value = "  +1.24%  "
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("E51").Value = "  +4.57%  "
